# Apply cryptos.xlsx price/volume updates (+ row 29/30 coin swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.846.98"
$ws.Range("E2").Value = "  +2.15%  "

$ws.Range("D3").Value = "'2.390.21"
$ws.Range("E3").Value = "  +1.49%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'555.72"
$ws.Range("E5").Value = "  +2.91%  "

$ws.Range("D6").Value = "'141.86"
$ws.Range("E6").Value = "  +4.20%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "'0.527"
$ws.Range("E8").Value = "  +1.05%  "

$ws.Range("D9").Value = "'2.391.42"
$ws.Range("E9").Value = "  +1.53%  "

$ws.Range("D10").Value = "'0.110"
$ws.Range("E10").Value = "  +5.62%  "

$ws.Range("E11").Value = "  +2.00%  "

$ws.Range("D12").Value = "'5.38"
$ws.Range("E12").Value = "  +2.72%  "

$ws.Range("E13").Value = "  +4.20%  "

$ws.Range("D14").Value = "'25.83"
$ws.Range("E14").Value = "  +5.37%  "

$ws.Range("D15").Value = "'0.0000176"
$ws.Range("E15").Value = "  +9.29%  "

$ws.Range("D16").Value = "'2.816.08"
$ws.Range("E16").Value = "  +1.45%  "

$ws.Range("D17").Value = "'61.553.62"
$ws.Range("E17").Value = "  +1.38%  "

$ws.Range("D18").Value = "'2.391.68"
$ws.Range("E18").Value = "  +1.57%  "

$ws.Range("D19").Value = "'11.05"
$ws.Range("E19").Value = "  +4.45%  "

$ws.Range("E20").Value = "  +3.44%  "

$ws.Range("D21").Value = "'322.87"
$ws.Range("E21").Value = "  +2.85%  "

$ws.Range("D22").Value = "'6.73"
$ws.Range("E22").Value = "  +2.28%  "

$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("D24").Value = "'1.78"
$ws.Range("E24").Value = "  -4.57%  "

$ws.Range("D25").Value = "'64.48"
$ws.Range("E25").Value = "  +2.44%  "

$ws.Range("D26").Value = "'8.91"
$ws.Range("E26").Value = "  +5.42%  "

$ws.Range("D27").Value = "'0.994"
$ws.Range("E27").Value = "  -0.72%  "

$ws.Range("D28").Value = "'2.501.58"
$ws.Range("E28").Value = "  +1.39%  "

$ws.Range("B29").Value = "Bittensor"
$ws.Range("C29").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D29").Value = "'531.51"
$ws.Range("E29").Value = "  +7.47%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'8.28"
$ws.Range("E30").Value = "  +5.03%  "

$ws.Range("D31").Value = "'0.0₃0918"
$ws.Range("E31").Value = "  +3.78%  "

$ws.Range("D32").Value = "'1.42"
$ws.Range("E32").Value = "  +3.48%  "

$ws.Range("E33").Value = "  +3.72%  "

$ws.Range("D34").Value = "'1.85"
$ws.Range("E34").Value = "  +4.36%  "

$ws.Range("D35").Value = "'1.53"
$ws.Range("E35").Value = "  +0.99%  "

$ws.Range("D36").Value = "'0.997"
$ws.Range("E36").Value = "  -0.29%  "

$ws.Range("D37").Value = "'5.66"
$ws.Range("E37").Value = "  +9.07%  "

$ws.Range("D38").Value = "'4.77"
$ws.Range("E38").Value = "  +5.62%  "

$ws.Range("D39").Value = "'1.91"
$ws.Range("E39").Value = "  +8.33%  "

$ws.Range("D40").Value = "'0.381"
$ws.Range("E40").Value = "  +2.53%  "

$ws.Range("D41").Value = "'18.62"
$ws.Range("E41").Value = "  +1.82%  "

$ws.Range("D42").Value = "'146.05"
$ws.Range("E42").Value = "  +5.60%  "

$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("D44").Value = "'41.52"
$ws.Range("E44").Value = "  +3.63%  "

$ws.Range("D45").Value = "'149.85"
$ws.Range("E45").Value = "  +6.19%  "

$ws.Range("D46").Value = "'2.21"
$ws.Range("E46").Value = "  +6.42%  "

$ws.Range("D47").Value = "'3.63"
$ws.Range("E47").Value = "  +3.69%  "

$ws.Range("D48").Value = "'0.0529"
$ws.Range("E48").Value = "  +4.36%  "

$ws.Range("D49").Value = "'20.12"
$ws.Range("E49").Value = "  +4.18%  "

$ws.Range("D50").Value = "'0.586"
$ws.Range("E50").Value = "  +3.21%  "

$ws.Range("D51").Value = "'0.0909"
$ws.Range("E51").Value = "  +1.35%  "
